$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cl = $m.CustomLayouts.Item(1)
$cs = $cl.ThemeColorScheme
Write-Host "CustomLayout ThemeColorScheme:" $cs
Write-Host "Count:" $cs.Count
